$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 41.59528396220897
$ws.Range("C2").Value = 50.30575324618277
$ws.Range("D2").Value = 98.75011268317969
$ws.Range("E2").Value = 99.00091646870558
$ws.Range("F2").Value = 98.45669512320343
$ws.Range("G2").Value = 97.47332459209602
$ws.Range("H2").Value = 96.17691562397441

$ws.Range("B3").Value = 37.53297727120641
$ws.Range("C3").Value = 50.3563619084732
$ws.Range("D3").Value = 99.70492488693743
$ws.Range("E3").Value = 99.10033583658607
$ws.Range("F3").Value = 98.70730575517618
$ws.Range("G3").Value = 97.69881161670166
$ws.Range("H3").Value = 96.29163647640291

$ws.Range("B4").Value = 41.34724418994332
$ws.Range("C4").Value = 50.51563318690332
$ws.Range("D4").Value = 98.5950986662085
$ws.Range("E4").Value = 98.76124790183353
$ws.Range("F4").Value = 98.49453169309061
$ws.Range("G4").Value = 97.56225160802107
$ws.Range("H4").Value = 96.10649638695197

$ws.Range("B5").Value = 42.36680814045455
$ws.Range("C5").Value = 50.06775740033259
$ws.Range("D5").Value = 98.58365928125019
$ws.Range("E5").Value = 98.86678555203926
$ws.Range("F5").Value = 98.43880222383864
$ws.Range("G5").Value = 97.53384616161114
$ws.Range("H5").Value = 96.29565939507032

$ws.Range("B6").Value = 43.62511358165602
$ws.Range("C6").Value = 50.30878547958765
$ws.Range("D6").Value = 98.7945287793539
$ws.Range("E6").Value = 98.90840127624824
$ws.Range("F6").Value = 98.40678341834372
$ws.Range("G6").Value = 97.69090155514478
$ws.Range("H6").Value = 96.12524327902642
